# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K" = strikeouts) for rows 2-37 with freshly regenerated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 9
    3  = 9
    4  = 9
    5  = 6
    6  = 10
    7  = 2
    8  = 6
    9  = 7
    10 = 9
    11 = 12
    12 = 8
    13 = 8
    14 = 6
    15 = 8
    16 = 3
    17 = 8
    18 = 8
    19 = 7
    20 = 3
    21 = 13
    22 = 5
    23 = 7
    24 = 7
    25 = 6
    26 = 5
    27 = 7
    28 = 4
    29 = 7
    30 = 6
    31 = 7
    32 = 2
    33 = 5
    34 = 3
    35 = 3
    36 = 7
    37 = 3
}

foreach ($row in ($kValues.Keys | Sort-Object)) {
    $ws.Range("G$row").Value = $kValues[$row]
}
